$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.883.90'
$ws.Range("E2").Value = '  -2.42%  '
$ws.Range("D3").Value = '3.458.45'
$ws.Range("E3").Value = '  -1.69%  '
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = "'586.33"
$ws.Range("E5").Value = '  -3.48%  '
$ws.Range("D6").Value = "'136.90"
$ws.Range("E6").Value = '  -4.40%  '
$ws.Range("D7").Value = '3.456.03'
$ws.Range("E7").Value = '  -1.63%  '
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").Value = "'0.487"
$ws.Range("E9").Value = '  -4.74%  '
$ws.Range("D10").Value = "'0.122"
$ws.Range("E10").Value = '  -6.41%  '
$ws.Range("D11").Value = "'7.09"
$ws.Range("E11").Value = '  -8.01%  '
$ws.Range("E12").Value = '  -7.12%  '
$ws.Range("D13").Value = '4.045.31'
$ws.Range("E13").Value = '  -1.77%  '
$ws.Range("E14").Value = '  -6.69%  '
$ws.Range("D15").Value = "'26.46"
$ws.Range("E15").Value = '  -7.69%  '
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '3.459.13'
$ws.Range("E16").Value = '  -1.63%  '
$ws.Range("B17").Value = 'TRON'
$ws.Range("C17").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D17").Value = "'0.115"
$ws.Range("E17").Value = '  -1.52%  '
$ws.Range("D18").Value = '64.852.80'
$ws.Range("E18").Value = '  -2.25%  '
$ws.Range("D19").Value = "'9.62"
$ws.Range("E19").Value = '  -10.78%  '
$ws.Range("D20").Value = "'5.71"
$ws.Range("E20").Value = '  -7.30%  '
$ws.Range("D21").Value = "'13.70"
$ws.Range("E21").Value = '  -6.43%  '
$ws.Range("D22").Value = "'385.20"
$ws.Range("E22").Value = '  -8.76%  '
$ws.Range("D23").Value = "'0.551"
$ws.Range("E23").Value = '  -6.34%  '
$ws.Range("E24").Value = '  -0.09%  '
$ws.Range("D25").Value = "'72.25"
$ws.Range("E25").Value = '  -6.17%  '
$ws.Range("D26").Value = "'5.74"
$ws.Range("E26").Value = '  -0.24%  '
$ws.Range("D27").Value = '3.596.06'
$ws.Range("E27").Value = '  -1.79%  '
$ws.Range("D28").Value = "'0.0000107"
$ws.Range("E28").Value = '  -5.79%  '
$ws.Range("D29").Value = "'0.997"
$ws.Range("E29").Value = '  -0.36%  '
$ws.Range("D30").Value = "'7.29"
$ws.Range("E30").Value = '  -8.09%  '
$ws.Range("D31").Value = "'8.12"
$ws.Range("E31").Value = '  -8.90%  '
$ws.Range("D32").Value = "'2.20"
$ws.Range("E32").Value = '  -10.65%  '
$ws.Range("D33").Value = '3.474.89'
$ws.Range("E33").Value = '  -1.51%  '
$ws.Range("E34").Value = '  -0.02%  '
$ws.Range("D35").Value = "'22.94"
$ws.Range("E35").Value = '  -5.12%  '
$ws.Range("D36").Value = "'0.142"
$ws.Range("E36").Value = '  -8.30%  '
$ws.Range("D37").Value = "'170.00"
$ws.Range("E37").Value = '  -2.33%  '
$ws.Range("E38").Value = '  -10.89%  '
$ws.Range("D39").Value = "'6.81"
$ws.Range("E39").Value = '  -9.60%  '
$ws.Range("D40").Value = "'1.48"
$ws.Range("E40").Value = '  -9.24%  '
$ws.Range("D41").Value = "'4.69"
$ws.Range("E41").Value = '  -9.53%  '
$ws.Range("E42").Value = '  -6.15%  '
$ws.Range("E43").Value = '  -5.13%  '
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = '  +0.10%  '
$ws.Range("D45").Value = "'42.03"
$ws.Range("E45").Value = '  -7.25%  '
$ws.Range("E46").Value = '  -13.15%  '
$ws.Range("D47").Value = "'1.61"
$ws.Range("E47").Value = '  -8.28%  '
$ws.Range("D48").Value = "'23.26"
$ws.Range("E48").Value = '  +2.19%  '
$ws.Range("D49").Value = "'1.12"
$ws.Range("E49").Value = '  +1.48%  '
$ws.Range("E50").Value = '  -6.80%  '
$ws.Range("D51").Value = '2.216.63'
$ws.Range("E51").Value = '  -4.04%  '
